$wb = $excel.ActiveWorkbook

# --- 1. "unit" sheet: add three new rows (NR/Number, HAB_P/Inhabitants per ...,
#        P_HTHAB/Per hundred thousand inhabitants) --------------------------
$unit = $wb.Worksheets.Item("unit")

# Match the shared-string insertion order seen in the target workbook:
# names (column B) are written before codes (column A).
$unit.Range("B4").Value = "Number"
$unit.Range("B5").Value = "Inhabitants per ..."
$unit.Range("B6").Value = "Per hundred thousand inhabitants"
$unit.Range("A4").Value = "NR"
$unit.Range("A5").Value = "HAB_P"
$unit.Range("A6").Value = "P_HTHAB"

# Column A on this sheet carries the "code" style (style index 1 in the
# original file) -- copy it from an existing styled cell so no new style
# gets created.
$unit.Range("A2").Copy() | Out-Null
$unit.Range("A4:A6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- 2. new "victim" sheet, placed after "unit" -----------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$victim = $wb.Worksheets.Add($null, $lastSheet)
$victim.Name = "victim"

$victim.Range("A1").Value = "code"
$victim.Range("B1").Value = "name"
$victim.Range("A2").Value = "KIL"
$victim.Range("A3").Value = "INJ"
$victim.Range("B2").Value = "Killed"
$victim.Range("B3").Value = "Injured"

# Row 4 stays empty but column A keeps the "code" style, same as above.
$unit.Range("A2").Copy() | Out-Null
$victim.Range("A4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- 3. selections -----------------------------------------------------------
# Select on "unit" first, then on "victim" last so "victim" ends up the
# active / selected sheet.
$unit.Range("F11").Select()
$victim.Range("B1").Select()
